$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (this shifts the former rows 3..22 down to 4..23,
# and grows the used range from A1:R22 to A1:R23).
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new weekly data point.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = "10/13/2022"
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112012
$ws.Cells.Item(3, 7).Value = "Espinaca"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 6500
$ws.Cells.Item(3, 12).Value = 7000
$ws.Cells.Item(3, 13).Value = 6750
$ws.Cells.Item(3, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(3, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(3, 16).Value = 675
$ws.Cells.Item(3, 17).Value = 10
$ws.Cells.Item(3, 18).Value = "Hortaliza"
